$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4288.62
$ws.Range("I62").Value = 3861.3257
$ws.Range("K62").Value = 3861.3257
$ws.Range("M62").Value = -3237.3257
$ws.Range("H65").Value = 4288.62
$ws.Range("I65").Value = 3861.3257
$ws.Range("K65").Value = 19306.6285
$ws.Range("M65").Value = -16186.6285
$ws.Range("H69").Value = 379000
$ws.Range("I69").Value = 2000
$ws.Range("J69").Value = 504666.66
$ws.Range("K69").Value = 6000
$ws.Range("L69").Value = 1513999.98
$ws.Range("M69").Value = -5126
$ws.Range("N69").Value = -1515747.98
$ws.Range("H70").Value = 7008.095
$ws.Range("I70").Value = 6291.5
$ws.Range("J70").Value = 7963.5557
$ws.Range("K70").Value = 18874.5
$ws.Range("L70").Value = 23890.6671
$ws.Range("M70").Value = -18604.5
$ws.Range("N70").Value = -24430.6671
$ws.Range("H72").Value = 379000
$ws.Range("I72").Value = 2000
$ws.Range("J72").Value = 504666.66
$ws.Range("K72").Value = 18000
$ws.Range("L72").Value = 4541999.939999999
$ws.Range("M72").Value = -13632
$ws.Range("N72").Value = -4550735.939999999
$ws.Range("H73").Value = 7008.095
$ws.Range("I73").Value = 6291.5
$ws.Range("J73").Value = 7963.5557
$ws.Range("K73").Value = 18874.5
$ws.Range("L73").Value = 23890.6671
$ws.Range("M73").Value = -17938.5
$ws.Range("N73").Value = -25762.6671
$ws.Range("H98").Value = 1641.8649
$ws.Range("I98").Value = 1727.742
$ws.Range("J98").Value = 1198.1666
$ws.Range("K98").Value = 1727.742
$ws.Range("L98").Value = 1198.1666
$ws.Range("M98").Value = -229.742
$ws.Range("N98").Value = -4194.1666
$ws.Range("H116").Value = 16575.1
$ws.Range("I116").Value = 28553.75
$ws.Range("J116").Value = 8589.333000000001
$ws.Range("K116").Value = 28553.75
$ws.Range("L116").Value = 8589.333000000001
$ws.Range("M116").Value = -25111.75
$ws.Range("N116").Value = -15473.333
$ws.Range("H122").Value = 1641.8649
$ws.Range("I122").Value = 1727.742
$ws.Range("J122").Value = 1198.1666
$ws.Range("K122").Value = 5183.226
$ws.Range("L122").Value = 3594.4998
$ws.Range("M122").Value = -2733.226
$ws.Range("N122").Value = -8494.4998
$ws.Range("H136").Value = 125000
$ws.Range("J136").Value = 125000
$ws.Range("L136").Value = 125000
$ws.Range("N136").Value = -135200

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 7999
$ws.Range("I10").Value = 7999
$ws.Range("K10").Value = 7999
$ws.Range("M10").Value = -7829
$ws.Range("H32").Value = 16885.969
$ws.Range("I32").Value = 18047.934
$ws.Range("K32").Value = 18047.934
$ws.Range("M32").Value = -17760.934
$ws.Range("H61").Value = 4367.2896
$ws.Range("I61").Value = 783.1
$ws.Range("K61").Value = 783.1
$ws.Range("M61").Value = -571.1
$ws.Range("H63").Value = 3285.5715
$ws.Range("I63").Value = 2499.75
$ws.Range("K63").Value = 2499.75
$ws.Range("M63").Value = -1813.75
$ws.Range("H66").Value = 3285.5715
$ws.Range("I66").Value = 2499.75
$ws.Range("K66").Value = 12498.75
$ws.Range("M66").Value = -9066.75
$ws.Range("H74").Value = 211594
$ws.Range("I74").Value = 273531.38
$ws.Range("J74").Value = 16933.715
$ws.Range("K74").Value = 273531.38
$ws.Range("L74").Value = 16933.715
$ws.Range("M74").Value = -272657.38
$ws.Range("N74").Value = -18681.715
$ws.Range("H77").Value = 211594
$ws.Range("I77").Value = 273531.38
$ws.Range("J77").Value = 16933.715
$ws.Range("K77").Value = 1367656.9
$ws.Range("L77").Value = 84668.575
$ws.Range("M77").Value = -1363288.9
$ws.Range("N77").Value = -93404.575
$ws.Range("H88").Value = 5711
$ws.Range("I88").Value = 1689
$ws.Range("J88").Value = 7722
$ws.Range("K88").Value = 1689
$ws.Range("L88").Value = 7722
$ws.Range("M88").Value = -1283
$ws.Range("N88").Value = -8534
$ws.Range("H91").Value = 5711
$ws.Range("I91").Value = 1689
$ws.Range("J91").Value = 7722
$ws.Range("K91").Value = 1689
$ws.Range("L91").Value = 7722
$ws.Range("M91").Value = -285
$ws.Range("N91").Value = -10530
$ws.Range("H110").Value = 1448.75
$ws.Range("J110").Value = 2266.3333
$ws.Range("L110").Value = 2266.3333
$ws.Range("N110").Value = -6356.3333
$ws.Range("H136").Value = 4367.2896
$ws.Range("I136").Value = 783.1
$ws.Range("K136").Value = 2349.3
$ws.Range("M136").Value = 200.6999999999998

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 12855.481
$ws.Range("I20").Value = 18683.166
$ws.Range("J20").Value = 1200.1111
$ws.Range("K20").Value = 18683.166
$ws.Range("L20").Value = 1200.1111
$ws.Range("M20").Value = -18436.166
$ws.Range("N20").Value = -1694.1111
$ws.Range("H86").Value = 1513
$ws.Range("I86").Value = 1480.1538
$ws.Range("J86").Value = 1619.75
$ws.Range("K86").Value = 1480.1538
$ws.Range("L86").Value = 1619.75
$ws.Range("M86").Value = -357.1538
$ws.Range("N86").Value = -3865.75
$ws.Range("H89").Value = 1513
$ws.Range("I89").Value = 1480.1538
$ws.Range("J89").Value = 1619.75
$ws.Range("K89").Value = 7400.769
$ws.Range("L89").Value = 8098.75
$ws.Range("M89").Value = -1784.769
$ws.Range("N89").Value = -19330.75
$ws.Range("H94").Value = 2651.0625
$ws.Range("I94").Value = 1640
$ws.Range("J94").Value = 4336.1665
$ws.Range("K94").Value = 1640
$ws.Range("L94").Value = 4336.1665
$ws.Range("M94").Value = -1189
$ws.Range("N94").Value = -5238.1665
$ws.Range("H99").Value = 1728.4117
$ws.Range("I99").Value = 1424.1333
$ws.Range("K99").Value = 1424.1333
$ws.Range("M99").Value = 73.86670000000004

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").Value = $null
$ws.Range("H94").Value = 1716.1428
$ws.Range("I94").Value = 1468.125
$ws.Range("K94").Value = 1468.125
$ws.Range("M94").Value = -1017.125
$ws.Range("H99").Value = 5229.048
$ws.Range("I99").Value = 4021.5715
$ws.Range("J99").Value = 7644
$ws.Range("K99").Value = 4021.5715
$ws.Range("L99").Value = 7644
$ws.Range("M99").Value = -2523.5715
$ws.Range("N99").Value = -10640
$ws.Range("H107").Value = 1152.6
$ws.Range("J107").Value = 1152.6
$ws.Range("L107").Value = 1152.6
$ws.Range("N107").Value = -4992.6
$ws.Range("H122").Value = 1452.6296
$ws.Range("I122").Value = 1368.421
$ws.Range("J122").Value = 1652.625
$ws.Range("K122").Value = 4105.263
$ws.Range("L122").Value = 4957.875
$ws.Range("M122").Value = -1655.263
$ws.Range("N122").Value = -9857.875
$ws.Range("H126").Value = 5229.048
$ws.Range("I126").Value = 4021.5715
$ws.Range("J126").Value = 7644
$ws.Range("K126").Value = 12064.7145
$ws.Range("L126").Value = 22932
$ws.Range("M126").Value = -9594.7145
$ws.Range("N126").Value = -27872

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 835.62964
$ws.Range("J122").Value = 950.4
$ws.Range("L122").Value = 8553.6
$ws.Range("N122").Value = -13453.6
$ws.Range("H140").Value = 3012.0557
$ws.Range("I140").Value = 3012.0557
$ws.Range("K140").Value = 9036.167099999999
$ws.Range("M140").Value = -3856.167099999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1786.12
$ws.Range("I97").Value = 1139.9286
$ws.Range("K97").Value = 1139.9286
$ws.Range("M97").Value = -643.9286
$ws.Range("H126").Value = 3444.375
$ws.Range("I126").Value = 1447
$ws.Range("J126").Value = 4110.1665
$ws.Range("K126").Value = 4341
$ws.Range("L126").Value = 12330.4995
$ws.Range("M126").Value = -1871
$ws.Range("N126").Value = -17270.4995
$ws.Range("H132").Value = 1587.069
$ws.Range("I132").Value = 1549.04
$ws.Range("K132").Value = 4647.12
$ws.Range("M132").Value = -2117.12

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1284.1724
$ws.Range("I22").Value = 1186.9375
$ws.Range("K22").Value = 1186.9375
$ws.Range("M22").Value = -891.9375
$ws.Range("H27").Value = 1284.1724
$ws.Range("I27").Value = 1186.9375
$ws.Range("K27").Value = 1186.9375
$ws.Range("M27").Value = -1079.9375
$ws.Range("H86").Value = 79995
$ws.Range("J86").Value = 79995
$ws.Range("L86").Value = 79995
$ws.Range("N86").Value = -82367
$ws.Range("H89").Value = 79995
$ws.Range("J89").Value = 79995
$ws.Range("L89").Value = 239985
$ws.Range("N89").Value = -251841
$ws.Range("H93").Value = 3477
$ws.Range("I93").Value = 3138.6667
$ws.Range("J93").Value = 4492
$ws.Range("K93").Value = 3138.6667
$ws.Range("L93").Value = 4492
$ws.Range("M93").Value = -1890.6667
$ws.Range("N93").Value = -6988
$ws.Range("H100").Value = 4624.75
$ws.Range("I100").Value = 3000
$ws.Range("K100").Value = 3000
$ws.Range("M100").Value = -2459

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 10666
$ws.Range("J7").Value = 10666
$ws.Range("L7").Value = 10666
$ws.Range("N7").Value = -10892
$ws.Range("H14").Value = 3506.2856
$ws.Range("I14").Value = 3000
$ws.Range("J14").Value = 4329
$ws.Range("K14").Value = 3000
$ws.Range("L14").Value = 4329
$ws.Range("M14").Value = -2832
$ws.Range("N14").Value = -4665
